$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")

# Set B2:D24 to TRUE (was FALSE)
$ws.Range("B2:D24").Value = $true

# Set E2:E11 and E16:E24 to TRUE (was FALSE); E12:E15 remain FALSE (unchanged)
$ws.Range("E2:E11").Value = $true
$ws.Range("E16:E24").Value = $true

# Update the selected cell / active cell on the sheet
$ws.Activate()
$ws.Range("K17").Select()
